$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet (Date + Count) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-09-13T20:57:31+00:00"

# "Count" must stay a *text* value ("3"), not a number, to match the other
# Property/Value rows on this sheet. A leading apostrophe is the standard
# way to force a numeric-looking literal to be stored as text.
$meta.Range("B22").Value = "'3"

# --- Append a new concept row to the Concepts sheet ---
$concepts = $wb.Worksheets.Item("Concepts")

# Clone the formatting (and "Level" = 1 value) of the last existing data row
# (row 3) into the new row 4, then overwrite Code/Display with the new
# concept's data. Definition (D) stays blank, as with the other rows.
$concepts.Range("A3:D3").Copy($concepts.Range("A4:D4"))

$concepts.Range("B4").Value = "unknown"
$concepts.Range("C4").Value = "Unknown"
